$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("i18n")

# Insert 50 new rows before row 1092 (shifts old rows 1092-1101 down to 1142-1151)
$ws.Rows("1092:1141").Insert()

# Data table: key, list of (localeId, value)
$groups = @(
    @("solved", @(@(1, "Solved"), @(2, "தீர்க்கப்பட்டது"), @(3, "Résolu"), @(4, "解决了"), @(5, "解決済み"), @(6, "Resolvido"), @(7, "हल किया"), @(8, "Решено"), @(9, "Resuelto"), @(10, "تم حلها"))),
    @("ahead", @(@(1, "Ahead"), @(2, "முன்னாடி"), @(3, "en avant"), @(4, "先"), @(5, "先に"), @(6, "à frente"), @(7, "आगे"), @(8, "предстоящий"), @(9, "adelante"), @(10, "امام"))),
    @("tally", @(@(1, "Tally"), @(2, "எண்ணிக்கை"), @(3, "எண்ணிக்கை"), @(4, "相符"), @(5, "集計する"), @(6, "contar"), @(7, "गणना"), @(8, "подсчитывать"), @(9, "cuenta"), @(10, "حصيلة"))),
    @("unaccounted", @(@(1, "Unaccounted"), @(2, "கணக்கில் காட்டப்படாத"), @(3, "inexpliqué"), @(4, "下落不明"), @(5, "行方不明"), @(6, "não contabilizado"), @(7, "बेहिसाब"), @(8, "неучтенный"), @(9, "no contabilizado"), @(10, "في عداد المفقودين"))),
    @("wallet", @(@(1, "Wallet"), @(2, "பணப்பை"), @(3, "Porte monnaie"), @(4, "钱包"), @(5, "財布"), @(6, "Carteira"), @(7, "बटुआ"), @(8, "Бумажник"), @(9, "Cartera"), @(10, "محفظة"))),
)

$row = 1092
foreach ($group in $groups) {
    $key = $group[0]
    $items = $group[1]
    foreach ($item in $items) {
        $localeId = $item[0]
        $val = $item[1]
        $ws.Cells.Item($row, 2).Value = $localeId
        $ws.Cells.Item($row, 3).Value = $key
        $ws.Cells.Item($row, 4).Value = $val
        $row = $row + 1
    }
}

# Column A: running id formula (extends existing shared-formula pattern)
$ws.Range("A1092:A1141").Formula = "=A1091+1"

# Column E: CONCATENATE formula matching the existing pattern, extended over old+new range
$concatFormula = '=CONCATENATE("(",CHAR(34),A1081,CHAR(34),",",CHAR(34),B1081,CHAR(34),",",CHAR(34),C1081,CHAR(34),",",CHAR(34),D1081,CHAR(34),"),")'
$ws.Range("E1081:E1141").Formula = $concatFormula

# Restore workbook/view metadata changes from the diff
$ws.Application.ActiveWindow.ScrollRow = 1117
$ws.Range("D1144").Select()
